# Add three newly-completed "React" course entries to rows 18-20 of the
# "学習計画書" (study plan) sheet, and mark their progress (column G) as
# done (100%).
#
# Rows 18-20 previously had empty "実施内容" (F) and "進捗" (G) cells.
# We copy the formatting (incl. the mixed-run rich text used for the
# "React..." entries and the percentage style used for completed
# progress) from existing, already-filled rows, then overwrite the
# text/number so the shared-formatting (font, number format, borders)
# matches the rest of the table exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("学習計画書")

# F15 ("...React系列课程从零基础到项目开发实战 20-36") carries the same
# rich-text run formatting that the new F18:F20 entries need.
$reactTemplate = $ws.Cells.Item(15, 6)
# G4 is an existing "progress = 100%" cell (percentage number format).
$progressTemplate = $ws.Cells.Item(4, 7)

$rows = @(18, 19, 20)
$labels = @(
    "React系列课程从零基础到项目开发实战 37-49",
    "React系列课程从零基础到项目开发实战 49-60",
    "React系列课程从零基础到项目开发实战 61-63"
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]

    $fCell = $ws.Cells.Item($r, 6)
    $reactTemplate.Copy($fCell)
    $fCell.Value = $labels[$i]

    $gCell = $ws.Cells.Item($r, 7)
    $progressTemplate.Copy($gCell)
    $gCell.Value = 1
}

Write-Output "Updated F18:G20 with React course entries and progress"
